# PLANILLA DE EVALUACIÓN FASE 1 — "Planilla de evaluacion completada por la profe"
#
# The teacher finished grading: she re-ordered the three team members into
# their real roles (EVALUACION1!B4:B6), filled the two blank "X" cells in the
# individual IEP rubric (D13/D16, and their mirrored "L" column F13/F16) down
# from the already-completed rows above/below them, and then hand-typed a few
# final scores over what used to be live formulas (rows 20, 32, 44 and 55).
# Every downstream total/VLOOKUP/weighted-average cell is a formula already,
# so it recalculates on its own once the inputs above change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EVALUACION1")

# ---------------------------------------------------------------------
# 1) Re-order the three team members.
#    Was: B4 Alexander Hernández / B5 Álvaro Muñoz / B6 Abel Sánchez
#    Now: B4 Abel Sánchez / B5 "Álvaro Muñoz " (note trailing space) / B6 Alexander Hernández
# ---------------------------------------------------------------------
$ws.Range("B4").Value2 = "Abel Sánchez"
$ws.Range("B5").Value2 = "Álvaro Muñoz "
$ws.Range("B6").Value2 = "Alexander Hernández"

# ---------------------------------------------------------------------
# 2) Complete the individual-rubric "Nivel de Logro" table: fill the
#    CL ("X") and L ("X") formulas down through D13:D16 / F13:F16, matching
#    the pattern already present in rows 14-15.
# ---------------------------------------------------------------------
$ws.Range("D13").Formula = '=IF($C13=CL,"X","")'
$ws.Range("D16").Formula = '=IF($C16=CL,"X","")'

$ws.Range("F13").Formula = '=IF($C13=L,"X","")'
$ws.Range("F16").Formula = '=IF($C16=L,"X","")'

# ---------------------------------------------------------------------
# 3) Row 20 (9. Utiliza reglas de redacción ...): the teacher typed the
#    final score by hand instead of leaving the rubric formula in place.
# ---------------------------------------------------------------------
$ws.Range("E20").Value2 = 0
$ws.Range("I20").Value2 = 4

# ---------------------------------------------------------------------
# 4) Row 32 (persona 1 - "11. Expone el tema ..."): value kept but formula
#    dropped; G32 cleared out entirely.
# ---------------------------------------------------------------------
$ws.Range("E32").Value2 = 10
$ws.Range("G32").ClearContents()

# ---------------------------------------------------------------------
# 5) Row 44 (persona 2 - "11. Expone el tema ..."): hand-typed 9 (was 10).
# ---------------------------------------------------------------------
$ws.Range("E44").Value2 = 9

# ---------------------------------------------------------------------
# 6) Row 55 (persona 3 - "11. Expone el tema ..."): hand-typed 9 (was 10).
# ---------------------------------------------------------------------
$ws.Range("E55").Value2 = 9

# Everything else (C4:E6 weighted averages, C23/E23/I23/C24, C46/E46/C47,
# C57/E57/C58, C27/C39/C50 name mirrors, ...) is formula-driven and
# recalculates automatically.
$excel.CalculateFullRebuild()

# Keep the teacher's last-looked-at cell in view, like the saved workbook.
$ws.Range("E58").Select()
